$wb = $excel.ActiveWorkbook

$oldGuid = "462aa2ac-6895-49d2-a241-4d5f5395bc9b"
$newGuid = "79134831-0534-4f0a-988d-4df5b37a1c1c"

$oldHash = "9b95bb2cc16deec1711fdb1d7bfa784d3c50331d"
$newHash = "0babea501d2388536c3d82cc4294d287dfce8673"

$hyperlinkUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d2d89d247a0fd213acab4a055c4a19c629c65314/e2e/$oldGuid.md"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-29 00:56:34"

$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hyperlinkUrl, "", "", "e2e\$newGuid.md")

# --- zh-cn sheet ---
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-29 00:56:30"

$wsZhCn.Range("A2").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $hyperlinkUrl, "", "", "$newGuid.md")

# --- de-de sheet ---
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-29 00:56:34"

$wsDeDe.Range("A2").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $hyperlinkUrl, "", "", "$newGuid.md")
